$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Simplify the redline-rating / has-land / has-buildings labels and fix the
# RAIN CLT founding year (was stored as text "1980s", now a numeric 1980).
$ws.Range("H1").Value = "Has land"
$ws.Range("I1").Value = "Has buildings"
$ws.Range("K2").Value = "C"
$ws.Range("K6").Value = "D"
$ws.Range("G9").Value = 1980
$ws.Range("K11").Value = "D"
$ws.Range("K16").Value = "D"
$ws.Range("K17").Value = "D"

# Turn on AutoFilter for the data table and move the active selection.
$ws.Range("A1:N19").AutoFilter()
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$N`$19")
$n.Visible = $false

$ws.Range("I2").Select()
